$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New bug report row (row 4) ---------------------------------------
# Values are written in the same order the strings appear in the target
# sharedStrings table (Issue, Component, Description, Status) so the new
# shared-string entries land at the expected indices.
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Search crash in Edit Companies window"
$ws.Range("D4").Value = "EditCompaniesWindow"
$ws.Range("E4").Value = "Searching for a company that doesn’t exist will case a crash."
$ws.Range("C4").Value = "FIXED"
$ws.Range("F4").Value = Get-Date -Year 2010 -Month 3 -Day 8 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("G4").Value = Get-Date -Year 2010 -Month 3 -Day 8 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# --- Match the formatting already used by the table rows above --------
$ws.Range("A4:E4").WrapText = $true
$ws.Range("A4:E4").VerticalAlignment = -4160
$ws.Range("F4:G4").WrapText = $true
$ws.Range("F4:G4").VerticalAlignment = -4160
$ws.Range("F4:G4").NumberFormat = "yyyy\-mm\-dd;@"

# Row 4 is a two-line (wrapped) entry like row 3, so give it the same height
$ws.Rows.Item(4).RowHeight = 30

# The Component column needed to grow to fit the new text
$ws.Columns.Item(4).ColumnWidth = 24.6

# --- Selection mirrors the post-edit workbook state --------------------
$ws.Range("C4").Select()
